$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a second "import" row for the new bean (row 14)
$ws.Range("B14").Value = "import"
$ws.Range("C14").Value = "org.openl.rules.beans.B1"

# Add new method table "Method B1 hello1()" at rows 21-22
$ws.Range("B21").Value = "Method B1 hello1()"
$ws.Range("B22").Value = 'return B1(name="hello");'

# Add new method table "Method B1 hello2()" at rows 26-27
$ws.Range("B26").Value = "Method B1 hello2()"
$ws.Range("B27").Value = 'return B1(var="hello");'

# Update the selected cell to match the saved view state
$ws.Range("C9").Select()
